# Atualização de bases das ligas, do dia: 14-05-2024 às 01:09
#
# Changes applied to "Mexico Liga de Expansion":
#   1. Row 91 <-> Row 92 : all match data (columns B..AB) swapped between
#      the two rows (the ranking column A keeps its original 89 / 90).
#   2. Row 186 <-> Row 187 : same full data swap (column A keeps 184 / 185).
#   3. A brand-new match is appended as row 247 (Atlante vs Universidad
#      Guadalajara), extending the used range from A1:AB246 to A1:AB247.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2, $firstCol, $lastCol) {
    # Value2 round-trips plain scalars (numbers/strings) cleanly through
    # this COM shim, unlike Value which hands back an unevaluated
    # property accessor when captured into a variable.
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $val1 = $ws.Cells.Item($row1, $c).Value2
        $val2 = $ws.Cells.Item($row2, $c).Value2
        $ws.Cells.Item($row1, $c).Value = $val2
        $ws.Cells.Item($row2, $c).Value = $val1
    }
}

# 1) Swap match data between rows 91 and 92 (keep column A as-is).
Swap-RowData 91 92 2 28

# 2) Swap match data between rows 186 and 187 (keep column A as-is).
Swap-RowData 186 187 2 28

# 3) Append new row 247, cloning the formatting of the last existing
#    row (bold/bordered rank cell in A, custom date format in D) before
#    writing in the new match's values.
$ws.Range("A246:AB246").Copy()
$ws.Range("A247:AB247").PasteSpecial(-4122)

$ws.Cells.Item(247, 1).Value  = 245
$ws.Cells.Item(247, 2).Value  = 8185480
$ws.Cells.Item(247, 3).Value  = "Mexico Liga de Expansion"
$ws.Cells.Item(247, 4).Value  = 45420.91666666666
$ws.Cells.Item(247, 5).Value  = "Atlante"
$ws.Cells.Item(247, 6).Value  = "Universidad Guadalajara"
$ws.Cells.Item(247, 7).Value  = 2
$ws.Cells.Item(247, 8).Value  = 0
$ws.Cells.Item(247, 9).Value  = "H"
$ws.Cells.Item(247, 10).Value = 1.571
$ws.Cells.Item(247, 11).Value = 3.75
$ws.Cells.Item(247, 12).Value = 4.75
$ws.Cells.Item(247, 13).Value = 1.65
$ws.Cells.Item(247, 14).Value = 3.6
$ws.Cells.Item(247, 15).Value = 5.25
$ws.Cells.Item(247, 16).Value = -0.75
$ws.Cells.Item(247, 17).Value = 1.85
$ws.Cells.Item(247, 18).Value = 1.95
$ws.Cells.Item(247, 19).Value = 2.25
$ws.Cells.Item(247, 20).Value = 1.85
$ws.Cells.Item(247, 21).Value = 1.95
$ws.Cells.Item(247, 22).Value = 0.6499999999999999
$ws.Cells.Item(247, 23).Value = -1
$ws.Cells.Item(247, 24).Value = -1
$ws.Cells.Item(247, 25).Value = 0.8500000000000001
$ws.Cells.Item(247, 26).Value = -1
$ws.Cells.Item(247, 27).Value = -0.5
$ws.Cells.Item(247, 28).Value = 0.475
